$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Insert a new row before row 31 (shifts rows 31+ down by one)
$ws.Rows("31:31").Insert()

# Populate the new row 31 with the new default production user
$ws.Cells.Item(31, 1).Value = "Oleg_Babak"
$ws.Cells.Item(31, 2).Value = "Password1!"
$ws.Cells.Item(31, 4).Value = "CUSTOM_USER"
$ws.Cells.Item(31, 5).Value = "Smoke Test User"
$ws.Cells.Item(31, 6).Value = "N"
